$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values are swapped between the paired rows.
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

# Row pairs that swap their (A,B,D,E,F,G,H,Q,R) contents.
$pairs = @(
    @(23, 24),
    @(28, 29),
    @(30, 31)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"

        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2

        $ws.Range($addr1).Value2 = $v2
        $ws.Range($addr2).Value2 = $v1
    }
}
